$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Rename the "Context" header (column B) to "Tags"
$ws.Range("B1").Value = "Tags"

# Grow the "History" table from 6 to 10 columns (A1:F2 -> A1:J2)
$tbl.Resize($ws.Range("A1:J2"))

# Name the four new header cells
$ws.Range("G1").Value = "E1"
$ws.Range("H1").Value = "E2"
$ws.Range("I1").Value = "E3"
$ws.Range("J1").Value = "E4"

# Match the column widths Excel settled on after the table grew
$ws.Columns.Item(1).ColumnWidth = 22.4252
$ws.Columns.Item(2).ColumnWidth = 18.7509
$ws.Columns.Item(3).ColumnWidth = 33.7509
$ws.Columns.Item(4).ColumnWidth = 23.9255
$ws.Columns.Item(5).ColumnWidth = 15.5858
$ws.Columns.Item(6).ColumnWidth = 18.5858

# Leave the selection where the author left it
$ws.Range("I1").Select()
